$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "print"
$ws.Range("B3").Value = "error"

$ws.Range("B3").Select()
